$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gas Sorption Input")

# Top summary block (rows 2-4)
$ws.Range("H2").Value = 92.36860169715462
$ws.Range("I2").Value = 1.819015083414506

$ws.Range("H3").Value = 5.368215497467331
$ws.Range("I3").Value = 0.36648016077864065

$ws.Range("H4").Value = 15.381177201712239
$ws.Range("I4").Value = 0.4761645837834983

# Dual Mode prediction table (rows 12-18)
$ws.Range("P12").Value = 52.87600173852176
$ws.Range("Q12").Value = 1.8498667314869128

$ws.Range("P13").Value = 83.70692705335281
$ws.Range("Q13").Value = 1.8086625703807682

$ws.Range("P14").Value = 97.54526414383393
$ws.Range("Q14").Value = 1.8245274789006884

$ws.Range("P15").Value = 109.3122802817562
$ws.Range("Q15").Value = 1.9117807281521277

$ws.Range("P16").Value = 119.23485564423626
$ws.Range("Q16").Value = 2.03215487278341

$ws.Range("P17").Value = 128.32498050774956
$ws.Range("Q17").Value = 2.173085158612707

$ws.Range("P18").Value = 135.53970164476718
$ws.Range("Q18").Value = 2.3021739961767254
